$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly rows to append below the existing data (rows 22-24)
$data = @(
    @(45901, 0.14306, 0.21049, 0.54474, 0.10171, 0.05498),
    @(45931, 0.13266, 0.20834, 0.54833, 0.11067, 0.04024),
    @(45962, 0.13572, 0.18174, 0.48085, 0.20169, 0.04024)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
